$wb = $excel.ActiveWorkbook

# Sheet "展览" (Worksheet 1): update F4, F5, F6
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F4").Value = 2070
$ws1.Range("F5").Value = 163
$ws1.Range("F6").Value = 357

# Sheet "全部类型" (Worksheet 4): update F4, F5, F7
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value = 2070
$ws4.Range("F5").Value = 163
$ws4.Range("F7").Value = 357
